# CS107-68: Adds LCS and NW Algorithm similarity% in front-end
#
# studentData.xlsx test fixture - replace the last two placeholder rows
# ("h"/"i") with real student entries (Max / Kat) and give each of them
# their own mailto hyperlink instead of sharing the generic C4:C12 link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: h -> Max / max@x.com
$ws.Range("B11").Value = "Max"
$ws.Range("C11").Value = "max@x.com"

# Row 12: i -> Kat / kat@x.com
$ws.Range("B12").Value = "Kat"
$ws.Range("C12").Value = "kat@x.com"

# Give C11 / C12 their own hyperlinks (on top of the existing C4:C12 one),
# then restore the standard "Hyperlink" cell style so it matches the rest
# of the column (Hyperlinks.Add() otherwise stamps its own style variant).
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:max@x.com") | Out-Null
$ws.Range("C11").Style = $ws.Range("C2").Style

$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:kat@x.com") | Out-Null
$ws.Range("C12").Style = $ws.Range("C2").Style

# Leave the selection where the user ended up after typing the new rows.
$ws.Range("B13").Select() | Out-Null
